# chore: update Sheets via scheduled runner
# Refresh market-price-derived columns (H:N) for the affected leve rows
# across the ARM, BSM, CRP, LTW and WVR sheets.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 191313.77
$ws.Range("I74").Value = 223351.62
$ws.Range("J74").Value = 80413.46000000001
$ws.Range("K74").Value = 223351.62
$ws.Range("L74").Value = 80413.46000000001
$ws.Range("M74").Value = -222477.62
$ws.Range("N74").Value = -82161.46000000001

$ws.Range("H77").Value = 191313.77
$ws.Range("I77").Value = 223351.62
$ws.Range("J77").Value = 80413.46000000001
$ws.Range("K77").Value = 1116758.1
$ws.Range("L77").Value = 402067.3
$ws.Range("M77").Value = -1112390.1
$ws.Range("N77").Value = -410803.3

$ws.Range("H97").Value = 209.31818
$ws.Range("I97").Value = 211.19048
$ws.Range("J97").Value = 170
$ws.Range("K97").Value = 211.19048
$ws.Range("L97").Value = 170
$ws.Range("M97").Value = 284.80952
$ws.Range("N97").Value = -1162

$ws.Range("H110").Value = 938.875
$ws.Range("I110").Value = 930.1429000000001
$ws.Range("J110").Value = 1000
$ws.Range("K110").Value = 930.1429000000001
$ws.Range("L110").Value = 1000
$ws.Range("M110").Value = 1114.8571
$ws.Range("N110").Value = -5090

$ws.Range("H121").Value = 33235
$ws.Range("I121").Value = 0
$ws.Range("J121").Value = 33235
$ws.Range("K121").Value = 0
$ws.Range("L121").Value = 33235
$ws.Range("N121").Value = -36729

$ws.Range("H122").Value = 3075.7073
$ws.Range("I122").Value = 2685.28
$ws.Range("J122").Value = 3685.75
$ws.Range("K122").Value = 8055.84
$ws.Range("L122").Value = 11057.25
$ws.Range("M122").Value = -5605.84
$ws.Range("N122").Value = -15957.25

$ws.Range("H123").Value = 23999.5
$ws.Range("I123").Value = 0
$ws.Range("J123").Value = 23999.5
$ws.Range("K123").Value = 0
$ws.Range("L123").Value = 23999.5
$ws.Range("N123").Value = -33799.5

$ws.Range("H124").Value = 10000
$ws.Range("I124").Value = 0
$ws.Range("J124").Value = 10000
$ws.Range("K124").Value = 0
$ws.Range("L124").Value = 10000
$ws.Range("N124").Value = -19820

$ws.Range("H125").Value = 30715
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 30715
$ws.Range("K125").Value = 0
$ws.Range("L125").Value = 30715
$ws.Range("N125").Value = -40555

$ws.Range("H126").Value = 0
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 0

$ws.Range("H127").Value = 0
$ws.Range("I127").Value = 0
$ws.Range("J127").Value = 0
$ws.Range("K127").Value = 0
$ws.Range("L127").Value = 0

$ws.Range("H128").Value = 35500
$ws.Range("I128").Value = 0
$ws.Range("J128").Value = 35500
$ws.Range("K128").Value = 0
$ws.Range("L128").Value = 35500
$ws.Range("N128").Value = -45460

$ws.Range("H129").Value = 40000
$ws.Range("I129").Value = 0
$ws.Range("J129").Value = 40000
$ws.Range("K129").Value = 0
$ws.Range("L129").Value = 40000
$ws.Range("N129").Value = -50000

$ws.Range("H130").Value = 56582.25
$ws.Range("I130").Value = 0
$ws.Range("J130").Value = 56582.25
$ws.Range("K130").Value = 0
$ws.Range("L130").Value = 56582.25
$ws.Range("N130").Value = -66622.25

$ws.Range("H131").Value = 51615.8
$ws.Range("I131").Value = 0
$ws.Range("J131").Value = 51615.8
$ws.Range("K131").Value = 0
$ws.Range("L131").Value = 51615.8
$ws.Range("N131").Value = -61695.8

$ws.Range("H132").Value = 26495.023
$ws.Range("I132").Value = 38583.75
$ws.Range("J132").Value = 3929.4
$ws.Range("K132").Value = 115751.25
$ws.Range("L132").Value = 11788.2
$ws.Range("M132").Value = -113221.25
$ws.Range("N132").Value = -16848.2

$ws.Range("H133").Value = 38420.332
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 38420.332
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 38420.332
$ws.Range("N133").Value = -43480.332

$ws.Range("H134").Value = 28849.75
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 28849.75
$ws.Range("K134").Value = 0
$ws.Range("L134").Value = 28849.75
$ws.Range("N134").Value = -38989.75

$ws.Range("H135").Value = 20325
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 20325
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 20325
$ws.Range("N135").Value = -30465

$ws.Range("H137").Value = 30246.666
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 30246.666
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 30246.666
$ws.Range("N137").Value = -40446.666

$ws.Range("H138").Value = 23104.75
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 23104.75
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 23104.75
$ws.Range("N138").Value = -33384.75

$ws.Range("H139").Value = 39243.57
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 39243.57
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 39243.57
$ws.Range("N139").Value = -49523.57

$ws.Range("H140").Value = 40000
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 40000
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 40000
$ws.Range("N140").Value = -50360

$ws.Range("H141").Value = 34342.25
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 34342.25
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 34342.25
$ws.Range("N141").Value = -44702.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 7861.25
$ws.Range("I99").Value = 8698.571
$ws.Range("K99").Value = 8698.571
$ws.Range("M99").Value = -7200.571

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 616.6667
$ws.Range("I16").Value = 525
$ws.Range("J16").Value = 800
$ws.Range("K16").Value = 525
$ws.Range("L16").Value = 800
$ws.Range("M16").Value = -238
$ws.Range("N16").Value = -1374

$ws.Range("H31").Value = 3097.1628
$ws.Range("I31").Value = 2219.5557
$ws.Range("J31").Value = 4578.125
$ws.Range("K31").Value = 2219.5557
$ws.Range("L31").Value = 4578.125
$ws.Range("M31").Value = -1924.5557
$ws.Range("N31").Value = -5168.125

$ws.Range("H34").Value = 3097.1628
$ws.Range("I34").Value = 2219.5557
$ws.Range("J34").Value = 4578.125
$ws.Range("K34").Value = 2219.5557
$ws.Range("L34").Value = 4578.125
$ws.Range("M34").Value = -2017.5557
$ws.Range("N34").Value = -4982.125

$ws.Range("H58").Value = 5987.48
$ws.Range("I58").Value = 8260.643
$ws.Range("K58").Value = 8260.643
$ws.Range("M58").Value = -8057.643

$ws.Range("H99").Value = 52193.05
$ws.Range("I99").Value = 78325.08
$ws.Range("J99").Value = 3662.1428
$ws.Range("K99").Value = 78325.08
$ws.Range("L99").Value = 3662.1428
$ws.Range("M99").Value = -76827.08
$ws.Range("N99").Value = -6658.1428

$ws.Range("H113").Value = 616.6667
$ws.Range("I113").Value = 525
$ws.Range("J113").Value = 800
$ws.Range("K113").Value = 525
$ws.Range("L113").Value = 800
$ws.Range("M113").Value = 1645
$ws.Range("N113").Value = -5140

$ws.Range("H126").Value = 52193.05
$ws.Range("I126").Value = 78325.08
$ws.Range("J126").Value = 3662.1428
$ws.Range("K126").Value = 234975.24
$ws.Range("L126").Value = 10986.4284
$ws.Range("M126").Value = -232505.24
$ws.Range("N126").Value = -15926.4284

$ws.Range("H132").Value = 2287.75
$ws.Range("I132").Value = 1182.8
$ws.Range("J132").Value = 4129.3335
$ws.Range("K132").Value = 3548.4
$ws.Range("L132").Value = 12388.0005
$ws.Range("M132").Value = -1018.4
$ws.Range("N132").Value = -17448.0005

$ws.Range("H134").Value = 1546.1351
$ws.Range("I134").Value = 967.14813
$ws.Range("J134").Value = 3109.4
$ws.Range("K134").Value = 2901.44439
$ws.Range("L134").Value = 9328.200000000001
$ws.Range("M134").Value = -366.4443900000001
$ws.Range("N134").Value = -14398.2

$ws.Range("H136").Value = 5987.48
$ws.Range("I136").Value = 8260.643
$ws.Range("K136").Value = 24781.929
$ws.Range("M136").Value = -22231.929

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 4110.033
$ws.Range("I132").Value = 3316.9167
$ws.Range("J132").Value = 4638.778
$ws.Range("K132").Value = 9950.750100000001
$ws.Range("L132").Value = 13916.334
$ws.Range("M132").Value = -7420.750100000001
$ws.Range("N132").Value = -18976.334

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H128").Value = 49000
$ws.Range("J128").Value = 49000
$ws.Range("L128").Value = 49000
$ws.Range("N128").Value = -58960
